$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.06093418598175
$ws.Range("B1").Value = 1.386205554008484
$ws.Range("C1").Value = 9.025010108947754
$ws.Range("D1").Value = 2.386106967926025
$ws.Range("E1").Value = 1.280614733695984
